$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 602; this shifts existing rows 602:643 down to 603:644
# and keeps every cell's formula/value intact (matches the diff: row N
# for N=603..644 after the edit equals row N-1 before the edit).
$ws.Rows.Item(602).Insert()

# Populate the freshly inserted row 602 with the new data point
# (2026/01/11, Sunday, hour 16 -> 19), written so that column A stays a
# plain text date string (matching the rest of the column) instead of
# being auto-converted to a date serial number.
$ws.Range("A602").NumberFormat = "@"
$ws.Range("A602").Value = "2026/01/11"
$ws.Range("B602").Value = "日"
$ws.Range("C602").Value = 16
$ws.Range("D602").Value = 19
$ws.Range("A602").ClearFormats()
